$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear the placeholder values in D2:G2 (they become blank cells and
# are dropped entirely from the sheet, matching the removed <c> elements)
$ws.Range("D2:G2").Value = ""

# Row 3: new entry -- "12/14/2024" kept as text (not converted to a date
# serial), "1" kept as text, remaining cells present but blank
$ws.Range("A3:H3").NumberFormat = "@"
$ws.Range("A3").Value = "12/14/2024"
$ws.Range("G3").Value = "1"

# Row 4: new entry -- only C4 has content
$ws.Range("A4:H4").NumberFormat = "@"
$ws.Range("C4").Value = "Missing Options!"

# Row 5: new entry -- G5 and H5 have content
$ws.Range("A5:H5").NumberFormat = "@"
$ws.Range("G5").Value = "wqew"
$ws.Range("H5").Value = "Freezes well,Freezes Poorly,Holds well in Fridge"
